$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Sheet1's selection moves off the previously-selected cell back to the
# top-left of the data (set first so the final Activate() below - done on
# Sheet2 - is what leaves the workbook's active tab on Sheet2).
$ws1.Range("A4").Select()

# --- New "Invoice -> Payment" mapping columns (H/I), mirroring the existing
#     Quote/SalesOrder/Shipping/Invoice column pairs already on the sheet ---
$ws2.Range("H1").Value = "Column Name"
$ws2.Range("I1").ClearFormats()
$ws2.Range("I1").Value = "TC003_Invoice_To_Payment"

$ws2.Range("H2").Value = "Payment  Number"
$ws2.Range("I2").Value = "PID-69"

$ws2.Range("H3").Value = "Payment  Status"
$ws2.Range("I3").Value = "PAID " + [char]0x2014 + " Partially Shipped"

# --- Refreshed quote / sales order / shipping / invoice reference numbers ---
$ws2.Range("B2").Value = "2021-90 " + [char]0x2014 + " ETCC UFT licence quote"
$ws2.Range("E2").Value = "2021-57 " + [char]0x2014 + " ETCC UFT licence Shipping"
$ws2.Range("B4").Value = "2021-63 " + [char]0x2014 + " ETCC UFT licence Sales Order"
$ws2.Range("E4").Value = "2021-77" + [char]10 + "ETCC UFT licence Invoice"

# Row 4 picked up an auto row height from the embedded line break above;
# re-fit it back down to the sheet's normal (default) row height.
$ws2.Cells.Item(4, 5).EntireRow.AutoFit()

# --- Column I needs to be wide enough for the new header text ---
$ws2.Columns.Item(9).ColumnWidth = 27

# --- Sheet2 becomes the active sheet/tab, scrolled to show the new columns,
#     with the new header cell selected ---
$ws2.Activate()
$ws2.Application.ActiveWindow.ScrollColumn = 3
$ws2.Range("H1").Select()
